$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.174.91"
$ws.Range("E2").Value = "  +0.89%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.854.57"
$ws.Range("E3").Value = "  +0.53%  "

$ws.Range("E4").Value = "  +1.00%  "

$ws.Range("E5").Value = "  +1.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "310.49"
$ws.Range("E6").Value = "  +0.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4774"
$ws.Range("E7").Value = "  +2.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3703"
$ws.Range("E8").Value = "  +1.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07281"
$ws.Range("E9").Value = "  +1.88%  "

$ws.Range("E10").Value = "  +0.96%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.96"
$ws.Range("E11").Value = "  +1.90%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07812"
$ws.Range("E12").Value = "  +1.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.851.46"
$ws.Range("E13").Value = "  -0.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.398"
$ws.Range("E14").Value = "  +2.18%  "

$ws.Range("E15").Value = "  +1.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.69"
$ws.Range("E16").Value = "  +1.51%  "

$ws.Range("E17").Value = "  +0.98%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008712"
$ws.Range("E18").Value = "  +1.05%  "

$ws.Range("E19").Value = "  +0.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.199.66"
$ws.Range("E20").Value = "  +0.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.65"
$ws.Range("E21").Value = "  +1.50%  "

$ws.Range("E22").Value = "  +1.25%  "

$ws.Range("E23").Value = "  +0.42%  "

$ws.Range("E24").Value = "  +1.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.16"
$ws.Range("E25").Value = "  +0.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.44"
$ws.Range("E26").Value = "  +1.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.995"
$ws.Range("E27").Value = "  -1.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.10"
$ws.Range("E28").Value = "  +0.67%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.931"
$ws.Range("E29").Value = "  +1.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08872"
$ws.Range("E30").Value = "  +0.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.308"
$ws.Range("E31").Value = "  +2.77%  "

$ws.Range("E32").Value = "  +0.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.558"
$ws.Range("E33").Value = "  +1.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7371"
$ws.Range("E34").Value = "  -1.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.689"
$ws.Range("E35").Value = "  -3.53%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.118"
$ws.Range("E36").Value = "  +3.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02000"
$ws.Range("E37").Value = "  +3.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05249"
$ws.Range("E38").Value = "  +1.23%  "

$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5301"
$ws.Range("E39").Value = "  +2.18%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.051"
$ws.Range("E40").Value = "  +2.10%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1529"
$ws.Range("E41").Value = "  +0.64%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.334"
$ws.Range("E42").Value = "  +2.47%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.64"
$ws.Range("E43").Value = "  +1.33%  "

$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4757"
$ws.Range("E44").Value = "  +1.45%  "

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.018"
$ws.Range("E45").Value = "  +1.18%  "

$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.38"
$ws.Range("E46").Value = "  +2.22%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.628"
$ws.Range("E47").Value = "  +1.41%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "66.12"
$ws.Range("E48").Value = "  +0.77%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06068"
$ws.Range("E49").Value = "  +0.38%  "

$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.8956"
$ws.Range("E50").Value = "  +0.45%  "

$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.58"
$ws.Range("E51").Value = "  +1.00%  "
